$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.731.60"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "2.669.27"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.77"
$ws.Range("E5").Value = "  -2.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.12"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.614"
$ws.Range("E8").Value = "  +3.98%  "
$ws.Range("E9").Value = "  +2.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.400"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("E11").Value = "  -3.83%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.07"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000198"
$ws.Range("E14").Value = "  -5.01%  "
$ws.Range("D15").Value = "3.148.26"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").Value = "65.583.87"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "2.663.04"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.68"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.79"
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.80"
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.49"
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.16"
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.67"
$ws.Range("E25").Value = "  -3.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.69"
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("E27").Value = "  -3.76%  "
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("E29").Value = "  -3.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "535.88"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  -4.08%  "
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.45"
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.48"
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "156.82"
$ws.Range("E39").Value = "  -3.42%  "
$ws.Range("E40").Value = "  -3.34%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "162.69"
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.32"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.62"
$ws.Range("E46").Value = "  -5.08%  "
$ws.Range("E47").Value = "  -4.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.638"
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("D49").Value = "0.0₆0256"
$ws.Range("E49").Value = "  +9.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.95"
$ws.Range("E50").Value = "  -4.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0987"
$ws.Range("E51").Value = "  -1.04%  "
